# Fixing the big mistake: correct Total (B) and Community (D) consumption values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2"  = 1616.739107633333
    "D2"  = 120.6574074666667
    "B3"  = 1507.6620386
    "D3"  = 111.7394087333333
    "B4"  = 1595.913612216666
    "D4"  = 114.86206215
    "B5"  = 1545.856946
    "D5"  = 115.4490732666667
    "B6"  = 1605.0634982
    "D6"  = 110.6163707333333
    "B7"  = 1547.514503649999
    "D7"  = 116.4257378333333
    "B8"  = 1608.401537199999
    "D8"  = 118.6550773333333
    "B9"  = 1598.758053183333
    "D9"  = 116.83039675
    "B10" = 1549.886868583333
    "D10" = 105.2393646833333
    "B11" = 1610.85472625
    "D11" = 117.1247235
    "B12" = 1551.372058616666
    "D12" = 123.3677795833333
    "B13" = 1548.880400366666
    "D13" = 116.8627440333333
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
